# Fix UI, Complete Data Permission
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update book titles (145 -> 245, 146 -> 246, 147 -> 247)
$ws.Range("A3").Value = "كتاب رقم 245 "
$ws.Range("A4").Value = "كتاب رقم 246"
$ws.Range("A5").Value = "كتاب رقم 247"

# Fix UI: move active selection from A9 to A7
$ws.Range("A7").Select()
